# Massive Changes (Added Printing and Semester)
# Applies the edits captured in the target diff:
#  - normalize a handful of Subj_timeout (column B) time values to their
#    canonical floating point representation (re-entering the same time)
#  - change row 16's "days" value from MONTHUFRI to MON
#  - fill in the previously-blank row 17 with a full data row (wed)
#  - append three more data rows (18: wed: 19/20: MON, with a room typo
#    "305 ICBs" on the last one)
#  - update the view: zoom 70%, scrolled/selected near the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize time values in column B (Subj_timeout) -----------------
# These re-assignments collapse the stored double to the same bit pattern
# Excel itself uses when the time is (re)entered, matching the target file.
$ws.Cells.Item(4, 2).Value  = 0.47916666666666669
$ws.Cells.Item(5, 2).Value  = 0.52083333333333337
$ws.Cells.Item(6, 2).Value  = 0.47916666666666669
$ws.Cells.Item(7, 2).Value  = 0.47916666666666669
$ws.Cells.Item(10, 2).Value = 0.39583333333333331
$ws.Cells.Item(12, 2).Value = 0.47916666666666669
$ws.Cells.Item(14, 2).Value = 0.39583333333333331
$ws.Cells.Item(16, 2).Value = 0.39583333333333331

# --- Row 16: "days" changes from MONTHUFRI to MON ----------------------
$ws.Cells.Item(16, 3).Value = "MON"

# --- Row 17: fill in the previously-empty row (keep its time style) ---
$ws.Cells.Item(17, 1).NumberFormat = $ws.Cells.Item(16, 1).NumberFormat
$ws.Cells.Item(17, 2).NumberFormat = $ws.Cells.Item(16, 2).NumberFormat
$ws.Cells.Item(17, 1).Value = 0.3125
$ws.Cells.Item(17, 2).Value = 0.39583333333333331
$ws.Cells.Item(17, 3).Value = "wed"
$ws.Cells.Item(17, 4).Value = "GEC001"
$ws.Cells.Item(17, 5).Value = "Purposive Communication"
$ws.Cells.Item(17, 6).Value = 3
$ws.Cells.Item(17, 7).Value = "305 ICB"
$ws.Cells.Item(17, 8).Value = "IT 1A"
$ws.Cells.Item(17, 9).Value = "ARISTOTLE"
$ws.Cells.Item(17, 10).Value = "CARANDANG"
$ws.Cells.Item(17, 11).Value = "F"
$ws.Cells.Item(17, 12).Value = "m"

# --- Row 18: same content as row 17 ------------------------------------
$ws.Cells.Item(18, 1).NumberFormat = $ws.Cells.Item(16, 1).NumberFormat
$ws.Cells.Item(18, 2).NumberFormat = $ws.Cells.Item(16, 2).NumberFormat
$ws.Cells.Item(18, 1).Value = 0.3125
$ws.Cells.Item(18, 2).Value = 0.39583333333333331
$ws.Cells.Item(18, 3).Value = "wed"
$ws.Cells.Item(18, 4).Value = "GEC001"
$ws.Cells.Item(18, 5).Value = "Purposive Communication"
$ws.Cells.Item(18, 6).Value = 3
$ws.Cells.Item(18, 7).Value = "305 ICB"
$ws.Cells.Item(18, 8).Value = "IT 1A"
$ws.Cells.Item(18, 9).Value = "ARISTOTLE"
$ws.Cells.Item(18, 10).Value = "CARANDANG"
$ws.Cells.Item(18, 11).Value = "F"
$ws.Cells.Item(18, 12).Value = "m"

# --- Row 19: MON / IT 2B -------------------------------------------------
$ws.Cells.Item(19, 1).NumberFormat = $ws.Cells.Item(16, 1).NumberFormat
$ws.Cells.Item(19, 2).NumberFormat = $ws.Cells.Item(16, 2).NumberFormat
$ws.Cells.Item(19, 1).Value = 0.3125
$ws.Cells.Item(19, 2).Value = 0.39583333333333331
$ws.Cells.Item(19, 3).Value = "MON"
$ws.Cells.Item(19, 4).Value = "GEC001"
$ws.Cells.Item(19, 5).Value = "Purposive Communication"
$ws.Cells.Item(19, 6).Value = 3
$ws.Cells.Item(19, 7).Value = "305 ICB"
$ws.Cells.Item(19, 8).Value = "IT 2B"
$ws.Cells.Item(19, 9).Value = "ARISTOTLE"
$ws.Cells.Item(19, 10).Value = "CARANDANG"
$ws.Cells.Item(19, 11).Value = "F"
$ws.Cells.Item(19, 12).Value = "m"

# --- Row 20: MON / IT 2B, with room "305 ICBs" ---------------------------
$ws.Cells.Item(20, 1).NumberFormat = $ws.Cells.Item(16, 1).NumberFormat
$ws.Cells.Item(20, 2).NumberFormat = $ws.Cells.Item(16, 2).NumberFormat
$ws.Cells.Item(20, 1).Value = 0.3125
$ws.Cells.Item(20, 2).Value = 0.39583333333333331
$ws.Cells.Item(20, 3).Value = "MON"
$ws.Cells.Item(20, 4).Value = "GEC001"
$ws.Cells.Item(20, 5).Value = "Purposive Communication"
$ws.Cells.Item(20, 6).Value = 3
$ws.Cells.Item(20, 7).Value = "305 ICBs"
$ws.Cells.Item(20, 8).Value = "IT 2B"
$ws.Cells.Item(20, 9).Value = "ARISTOTLE"
$ws.Cells.Item(20, 10).Value = "CARANDANG"
$ws.Cells.Item(20, 11).Value = "F"
$ws.Cells.Item(20, 12).Value = "m"

# --- View: zoom to 70%, scroll/select near the newly added rows --------
$excel.ActiveWindow.Zoom = 70
$null = $ws.Range("J20:L20").Select()
